$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2210796915167095
$ws.Range("C2").Value = 0.532133676092545
$ws.Range("J2").Value = 0.02570694087403599
$ws.Range("P2").Value = 0.141388174807198
$ws.Range("S2").Value = 0.07969151670951156
$ws.Range("C3").Value = 0.01428571428571429
$ws.Range("J3").Value = 0.03333333333333333
$ws.Range("P3").Value = 0.7571428571428571
$ws.Range("S3").Value = 0.1952380952380952
$ws.Range("J4").Value = 0.1186440677966102
$ws.Range("O4").Value = 0.01694915254237288
$ws.Range("P4").Value = 0.711864406779661
$ws.Range("S4").Value = 0.1525423728813559
$ws.Range("D6").Value = 0.01428571428571429
$ws.Range("F6").Value = 0.0761904761904762
$ws.Range("J6").Value = 0.280952380952381
$ws.Range("O6").Value = 0.01904761904761905
$ws.Range("Q6").Value = 0.1714285714285714
$ws.Range("R6").Value = 0.02857142857142857
$ws.Range("S6").Value = 0.3428571428571429
$ws.Range("B7").Value = 0.1179775280898876
$ws.Range("D7").Value = 0.02247191011235955
$ws.Range("F7").Value = 0.06179775280898876
$ws.Range("J7").Value = 0.2022471910112359
$ws.Range("O7").Value = 0.005617977528089887
$ws.Range("Q7").Value = 0.1404494382022472
$ws.Range("R7").Value = 0.08426966292134831
$ws.Range("S7").Value = 0.3651685393258427
$ws.Range("B8").Value = 0.1074168797953964
$ws.Range("D8").Value = 0.02557544757033248
$ws.Range("F8").Value = 0.04603580562659847
$ws.Range("J8").Value = 0.1687979539641944
$ws.Range("O8").Value = 0.01023017902813299
$ws.Range("Q8").Value = 0.1585677749360614
$ws.Range("R8").Value = 0.06649616368286446
$ws.Range("S8").Value = 0.4168797953964195
$ws.Range("B9").Value = 0.1171171171171171
$ws.Range("D9").Value = 0.02702702702702703
$ws.Range("E9").Value = 0.004504504504504504
$ws.Range("F9").Value = 0.07657657657657657
$ws.Range("J9").Value = 0.2432432432432433
$ws.Range("O9").Value = 0.009009009009009009
$ws.Range("Q9").Value = 0.1441441441441441
$ws.Range("R9").Value = 0.0945945945945946
$ws.Range("S9").Value = 0.2837837837837838
$ws.Range("B10").Value = 0.1335149863760218
$ws.Range("D10").Value = 0.02452316076294278
$ws.Range("F10").Value = 0.06198910081743869
$ws.Range("J10").Value = 0.1777929155313351
$ws.Range("O10").Value = 0.01226158038147139
$ws.Range("Q10").Value = 0.1934604904632153
$ws.Range("R10").Value = 0.07425068119891008
$ws.Range("S10").Value = 0.3222070844686649
$ws.Range("G11").Value = 0.1611721611721612
$ws.Range("J11").Value = 0.08424908424908426
$ws.Range("K11").Value = 0.2234432234432235
$ws.Range("L11").Value = 0.5128205128205128
$ws.Range("S11").Value = 0.01831501831501832
$ws.Range("F12").Value = 0.007246376811594203
$ws.Range("G12").Value = 0.7681159420289855
$ws.Range("J12").Value = 0.1594202898550725
$ws.Range("K12").Value = 0.01449275362318841
$ws.Range("L12").Value = 0.02173913043478261
$ws.Range("S12").Value = 0.02898550724637681
$ws.Range("G13").Value = 0.7173913043478261
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.02173913043478261
$ws.Range("F15").Value = 0.008510638297872341
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.07659574468085106
$ws.Range("J15").Value = 0.4212765957446808
$ws.Range("K15").Value = 0.05531914893617021
$ws.Range("M15").Value = 0.02553191489361702
$ws.Range("O15").Value = 0.05106382978723404
$ws.Range("S15").Value = 0.1617021276595745
$ws.Range("F16").Value = 0.016
$ws.Range("H16").Value = 0.144
$ws.Range("I16").Value = 0.112
$ws.Range("J16").Value = 0.428
$ws.Range("K16").Value = 0.08400000000000001
$ws.Range("M16").Value = 0.024
$ws.Range("N16").Value = 0.004
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.148
$ws.Range("F17").Value = 0.01363636363636364
$ws.Range("H17").Value = 0.1590909090909091
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.4340909090909091
$ws.Range("K17").Value = 0.07954545454545454
$ws.Range("M17").Value = 0.01363636363636364
$ws.Range("N17").Value = 0.002272727272727273
$ws.Range("O17").Value = 0.05454545454545454
$ws.Range("S17").Value = 0.1431818181818182
$ws.Range("H18").Value = 0.1371428571428571
$ws.Range("I18").Value = 0.08571428571428572
$ws.Range("J18").Value = 0.4571428571428571
$ws.Range("K18").Value = 0.1085714285714286
$ws.Range("M18").Value = 0.02857142857142857
$ws.Range("N18").Value = 0.005714285714285714
$ws.Range("O18").Value = 0.07428571428571429
$ws.Range("S18").Value = 0.1028571428571429
$ws.Range("F19").Value = 0.01937657961246841
$ws.Range("H19").Value = 0.1786015164279697
$ws.Range("I19").Value = 0.09688289806234204
$ws.Range("J19").Value = 0.3833192923336142
$ws.Range("K19").Value = 0.09772535804549284
$ws.Range("M19").Value = 0.01937657961246841
$ws.Range("N19").Value = 0.001684919966301601
$ws.Range("O19").Value = 0.09519797809604044
$ws.Range("S19").Value = 0.1078348778433024
